$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 (col I) and IF (col J), matching the header style
# already used by A1:H1 (bold font, thin box border, centered/top aligned)
# by copying H1's formatting onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in I0 and IF data for rows 2-68
$data = @(
    ,@(7,7)
    ,@(10,10)
    ,@(6,6)
    ,@(6,6)
    ,@(7,7)
    ,@(7,7)
    ,@(6,6)
    ,@(6,7)
    ,@(4,4)
    ,@(5,6)
    ,@(8,8)
    ,@(7,7)
    ,@(8,8)
    ,@(6,7)
    ,@(6,6)
    ,@(9,9)
    ,@(6,7)
    ,@(6,6)
    ,@(5,6)
    ,@(3,4)
    ,@(5,6)
    ,@(7,7)
    ,@(4,5)
    ,@(8,8)
    ,@(8,9)
    ,@(11,11)
    ,@(5,6)
    ,@(7,8)
    ,@(7,8)
    ,@(7,7)
    ,@(7,7)
    ,@(10,10)
    ,@(7,7)
    ,@(8,8)
    ,@(7,7)
    ,@(7,7)
    ,@(9,9)
    ,@(5,6)
    ,@(6,8)
    ,@(6,6)
    ,@(5,5)
    ,@(6,7)
    ,@(9,9)
    ,@(9,9)
    ,@(9,9)
    ,@(7,8)
    ,@(7,7)
    ,@(7,7)
    ,@(5,6)
    ,@(6,6)
    ,@(9,9)
    ,@(8,8)
    ,@(8,8)
    ,@(6,7)
    ,@(6,6)
    ,@(6,6)
    ,@(3,4)
    ,@(7,7)
    ,@(7,9)
    ,@(6,7)
    ,@(8,9)
    ,@(6,7)
    ,@(1,3)
    ,@(8,9)
    ,@(5,6)
    ,@(8,8)
    ,@(1,1)
)

for ($idx = 0; $idx -lt $data.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $data[$idx][0]
    $ws.Cells.Item($row, 10).Value = $data[$idx][1]
}

